$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update selected cell / view ---
$ws.Range("J3").Select()

# --- Row 5: add new values/formulas in D5:G5 ---
# Copy existing number formatting from row 7 (same layout) before writing values
$ws.Range("D7:F7").Copy()
$ws.Range("D5:F5").PasteSpecial(-4122)
$ws.Range("G7").Copy()
$ws.Range("G5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("D5").Value = 180
$ws.Range("E5").Value = 150
$ws.Range("F5").Formula = "=D5-E5"
$ws.Range("G5").Formula = "=E5/(E5+F5)"

# --- Row 7: update D7 and E7 ---
$ws.Range("D7").Value = 40
$ws.Range("E7").Value = 60

# --- Row 9: update E9, and change F9 from a formula to a plain value ---
$ws.Range("E9").Value = 26
$ws.Range("F9").Value = 0
